$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing SowDOY/HarvestDOY/EmDOY
# columns (G:I) to hold the new SowDate / HarvestDate / EmDate columns.
$ws.Columns("G:I").Insert()

# Headers for the three new columns.
$ws.Range("G1").Value = "SowDate"
$ws.Range("H1").Value = "HarvestDate"
$ws.Range("I1").Value = "EmDate"

# Sow date = 1 April, Harvest date = 1 November, Emergence date = 25 April
# for each trial year (rows 2-8 = years 2015-2021).
$sowDates = @("4/1/2015","4/1/2016","4/1/2017","4/1/2018","4/1/2019","4/1/2020","4/1/2021")
$harvestDates = @("11/1/2015","11/1/2016","11/1/2017","11/1/2018","11/1/2019","11/1/2020","11/1/2021")
$emDates = @("4/25/2015","4/25/2016","4/25/2017","4/25/2018","4/25/2019","4/25/2020","4/25/2021")

for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 2
    $ws.Range("G$r").Value = $sowDates[$i]
    $ws.Range("H$r").Value = $harvestDates[$i]
    $ws.Range("I$r").Value = $emDates[$i]
}

$ws.Range("G2:I8").NumberFormat = "mm-dd-yy"

$ws.Range("I4").Select()
